$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Atoms Office" contact row entirely.
$ws.Rows(9).Delete()

# Merge Country Code (col A) + Contact Number (col B) into a single WhatsApp
# number in column A, and move the contact's Name into column B.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $countryCode = $ws.Cells.Item($r, 1).Value2
    $contactNumber = $ws.Cells.Item($r, 2).Value2
    $name = $ws.Cells.Item($r, 3).Value2
    $combined = "$countryCode$contactNumber"
    $ws.Cells.Item($r, 1).Value = [double]$combined
    $ws.Cells.Item($r, 2).Value = $name
}

# New headers for the bulk-campaign CSV layout: number, then name.
# "name" is written first so it claims the earlier shared-string slot,
# matching how the sheet was actually authored.
$ws.Range("B1").Value = "name"
$ws.Range("A1").Value = "whatsappnumber"

# The old "Name" column (C) is no longer needed now that names live in column B.
$ws.Columns(3).ClearContents()

# Resize the two remaining columns to fit their new content
# (closest achievable best-fit widths for "whatsappnumber"/12-digit numbers
# and "name"/contact names).
$ws.Columns(1).ColumnWidth = 13.666666666666666
$ws.Columns(2).ColumnWidth = 10.833333333333332

$ws.Range("L19").Select()
